$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price/Volume columns store plain text values (numbers-as-text and
# percent strings like "-0.15%"). Force Text format before writing so
# Excel's COM layer does not auto-convert these into numeric/percentage
# values, then restore the original (default) style so no new cell
# formatting is introduced.
$dataRange = $ws.Range("D2:E50")
$origStyle = $dataRange.Style
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = "327.59"
$ws.Range("E2").Value = "-0.15%"
$ws.Range("E3").Value = "0.61%"
$ws.Range("D4").Value = "5.571"
$ws.Range("E4").Value = "-0.04%"
$ws.Range("D5").Value = "0.08062"
$ws.Range("E5").Value = "-0.53%"
$ws.Range("D6").Value = "1.965"
$ws.Range("E6").Value = "3.72%"
$ws.Range("D7").Value = "4.324"
$ws.Range("E7").Value = "1.02%"
$ws.Range("D8").Value = "2.554"
$ws.Range("E8").Value = "-7.31%"
$ws.Range("D9").Value = "0.9449"
$ws.Range("E9").Value = "0.85%"
$ws.Range("D10").Value = "0.1171"
$ws.Range("E10").Value = "-0.12%"
$ws.Range("D11").Value = "0.1864"
$ws.Range("E11").Value = "-1.64%"
$ws.Range("D12").Value = "11.81"
$ws.Range("E12").Value = "37.91%"
$ws.Range("D13").Value = "0.09865"
$ws.Range("E13").Value = "1.92%"
$ws.Range("D14").Value = "0.04746"
$ws.Range("E14").Value = "14.30%"
$ws.Range("D15").Value = "0.1065"
$ws.Range("E15").Value = "-0.32%"
$ws.Range("D16").Value = "0.001295"
$ws.Range("E16").Value = "1.40%"
$ws.Range("D17").Value = "0.04219"
$ws.Range("E17").Value = "-2.12%"
$ws.Range("D18").Value = "0.005886"
$ws.Range("E18").Value = "-0.93%"
$ws.Range("E19").Value = "-5.51%"
$ws.Range("E20").Value = "-0.30%"
$ws.Range("D21").Value = "0.1401"
$ws.Range("E21").Value = "2.72%"
$ws.Range("D22").Value = "0.2510"
$ws.Range("E22").Value = "-2.82%"
$ws.Range("D23").Value = "0.001256"
$ws.Range("E23").Value = "1.41%"
$ws.Range("D24").Value = "0.004298"
$ws.Range("E24").Value = "-2.20%"
$ws.Range("D25").Value = "0.0001192"
$ws.Range("E25").Value = "-3.01%"
$ws.Range("D26").Value = "0.0003752"
$ws.Range("E26").Value = "-5.84%"
$ws.Range("D38").Value = "0.02590"
$ws.Range("D39").Value = "0.05509"
$ws.Range("E39").Value = "0.54%"
$ws.Range("E40").Value = "-0.68%"
$ws.Range("D41").Value = "0.1400"
$ws.Range("E41").Value = "0.06%"
$ws.Range("D42").Value = "0.007487"
$ws.Range("E42").Value = "-34.51%"
$ws.Range("D43").Value = "0.002018"
$ws.Range("E43").Value = "-2.39%"
$ws.Range("D44").Value = "0.008354"
$ws.Range("E44").Value = "-13.28%"
$ws.Range("D45").Value = "0.00007096"
$ws.Range("E45").Value = "1.32%"
$ws.Range("E46").Value = "0.24%"
$ws.Range("E47").Value = "1.51%"
$ws.Range("D48").Value = "0.004837"
$ws.Range("E48").Value = "36.37%"
$ws.Range("D49").Value = "0.00002103"
$ws.Range("E49").Value = "0.24%"
$ws.Range("E50").Value = "0.24%"

$dataRange.Style = $origStyle
